# Add a new, blank "placeholder" data row (row 2) below the header row on
# the "Đơn sale chính" sheet, extending the sheet's used range from A1:T1
# to A1:T2.
#
# Text columns are set to an empty string, the numeric/money columns are
# set to 0, and "Mã dịch vụ" (column B) is left blank (it is a
# numeric-typed field with no value for this placeholder row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 2

# Text-valued columns (blank string):
# A  Tiền tố
# C  Ngày thực hiện
# D  Cơ sở
# E  Khách hàng
# F  Nguồn khách
# G  Tên dịch vụ
# H  Sale chính
# J  Sale phụ
# Q  Bác sĩ 1
# R  Bác sĩ 2
# S  Phụ phẫu 1
# T  Phụ phẫu 2
$textColumns = @("A", "C", "D", "E", "F", "G", "H", "J", "Q", "R", "S", "T")
foreach ($col in $textColumns) {
    $ws.Range("$col$row").Value = ""
}

# Numeric-valued columns (0):
# I  Đơn giá gốc
# K  Upsale
# L  Đơn giá
# M  Thanh toán lần đầu
# N  Trả sau
# O  Đã thanh toán
# P  Dư nợ
$numericColumns = @("I", "K", "L", "M", "N", "O", "P")
foreach ($col in $numericColumns) {
    $ws.Range("$col$row").Value = 0
}

# B  Mã dịch vụ - numeric field, left blank/empty for this placeholder row.
